$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 291, pushing existing rows 291:381 down to 292:382
$ws.Rows.Item(291).Insert()

# Populate the newly inserted row 291 with a new price-report record
$ws.Range("A291").Value = 4
$ws.Range("B291").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C291").Value = "Los Lagos"
$ws.Range("D291").Value = 44988
$ws.Range("E291").Value = 10
$ws.Range("F291").Value = "Fruta"
$ws.Range("G291").Value = 100108
$ws.Range("H291").Value = "Tropicales y subtropicales"
$ws.Range("I291").Value = 100108005
$ws.Range("J291").Value = "Piña"
$ws.Range("K291").Value = "Caramelo"
$ws.Range("L291").Value = "Primera"
$ws.Range("M291").Value = 200
$ws.Range("N291").Value = 25000
$ws.Range("O291").Value = 26000
$ws.Range("P291").Value = 25500
$ws.Range("Q291").Value = "$/caja 12 unidades"
$ws.Range("R291").Value = "Ecuador"
$ws.Range("S291").Value = 2125
$ws.Range("T291").Value = 12

# Match the date-formatted style used by the other "Fecha" column cells
$ws.Range("D291").NumberFormat = $ws.Range("D292").NumberFormat
